$wb = $excel.ActiveWorkbook

# --- Sheet "compare_models": update TT (Sec) column (I) for rows 2-19 ---
$wsCompare = $wb.Worksheets.Item("compare_models")

$wsCompare.Range("I2").Value = 0.074
$wsCompare.Range("I3").Value = 0.048
$wsCompare.Range("I4").Value = 0.042
$wsCompare.Range("I5").Value = 0.028
$wsCompare.Range("I6").Value = 0.092
$wsCompare.Range("I7").Value = 1.148
$wsCompare.Range("I9").Value = 0.018
$wsCompare.Range("I10").Value = 0.022
$wsCompare.Range("I11").Value = 0.02
$wsCompare.Range("I12").Value = 0.026
$wsCompare.Range("I13").Value = 0.594
$wsCompare.Range("I14").Value = 0.016
$wsCompare.Range("I15").Value = 0.022
$wsCompare.Range("I16").Value = 0.024
$wsCompare.Range("I17").Value = 0.022
$wsCompare.Range("I18").Value = 0.018
$wsCompare.Range("I19").Value = 0.026

# --- Sheet "pred_final": update metrics row 2 (C:H) ---
$wsPredFinal = $wb.Worksheets.Item("pred_final")

$wsPredFinal.Range("C2").Value = 1.6065
$wsPredFinal.Range("D2").Value = 4.2426
$wsPredFinal.Range("E2").Value = 2.0598
$wsPredFinal.Range("F2").Value = 0.9936
$wsPredFinal.Range("G2").Value = 0.0386
$wsPredFinal.Range("H2").Value = 0.0292
